$wb = $excel.ActiveWorkbook

# --- Work on the "optimization_parameters" sheet ---
$ws = $wb.Worksheets.Item("optimization_parameters")

# Insert a brand-new row above the current row 9 (shifts old rows 9-17 down to 10-18)
$ws.Rows.Item(9).Insert()

# Old row 8 held the "Model" label; rename it to "production_function" (value stays "Sigmoid")
$ws.Range("A8").Value = "production_function"

# The freshly inserted row 9 becomes the new "L_curve" parameter row
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 0

# Make this sheet the active / selected tab, with A9:B9 highlighted
$ws.Activate()
$ws.Range("A9:B9").Select()

$wb.Save()
